$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_R_acc_LG"

$values = @(
    96.839332748024589,
    97.67339771729587,
    97.67339771729587,
    97.9367866549605,
    97.980684811237921,
    98.024582967515357,
    97.278314310798947,
    97.585601404740999,
    98.419666374012294,
    98.331870061457423,
    96.312554872695344,
    98.639157155399474,
    96.663740122914845,
    97.585601404740999,
    97.058823529411768,
    97.410008779631255,
    98.46356453028973,
    97.14661984196664,
    96.619841966637395,
    97.410008779631255,
    97.89288849868305,
    96.971027216856882,
    93.942054433713778,
    94.337137840210715,
    98.244073748902551,
    97.322212467076383,
    97.585601404740999,
    97.497805092186127,
    97.805092186128178,
    98.770851624231781,
    98.46356453028973,
    98.507462686567166,
    98.419666374012294,
    98.112379280070243,
    97.717295873573306,
    75.021949078138718,
    97.14661984196664,
    98.024582967515357,
    96.312554872695344,
    98.419666374012294,
    96.883230904302025,
    97.278314310798947,
    98.244073748902551,
    97.366110623353819,
    97.805092186128178,
    97.234416154521512,
    85.952589991220378,
    98.244073748902551
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
